$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(7,6,2,2,6,8,5,5,3,3,2,7,1,5,3,5,4,4,4,4,2,4,3,3,7,4,2,3,3,0)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $kValues[$i]
}
